# "reprovando por falta corretamente"
#
# The sheet had been seeded with "_Modificado" suffixes tacked onto the
# header + every student's name (B3:B27). This restores the real names,
# and fills in the "Situação" (G) / "Nota para Aprovação Final" (H)
# columns, which were previously blank, using the attendance rule:
#   Faltas > 15  (i.e. missed more than 25% of the 60 total classes)
#       => "Reprovado por Falta"
#   else if average(P1,P2,P3) < 50
#       => "Reprovado por Nota"
#   else
#       => "Aprovado"
# "Nota para Aprovação Final" is 0 for every student in this pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Aluno_Modificado" -> "Aluno"
$ws.Cells.Item(3, 2).Value = "Aluno"

$totalAulas = 60

$alunos = @(
    @{ Row = 4;  Nome = "Eduardo" },
    @{ Row = 5;  Nome = "Murilo" },
    @{ Row = 6;  Nome = "Guilherme" },
    @{ Row = 7;  Nome = "Flavia " },
    @{ Row = 8;  Nome = "Ruan" },
    @{ Row = 9;  Nome = "Arnando" },
    @{ Row = 10; Nome = "Lucas" },
    @{ Row = 11; Nome = "Fabio" },
    @{ Row = 12; Nome = "Alisson" },
    @{ Row = 13; Nome = "Felipe" },
    @{ Row = 14; Nome = "Rachel" },
    @{ Row = 15; Nome = "Jouy" },
    @{ Row = 16; Nome = "François" },
    @{ Row = 17; Nome = "Dâmaris" },
    @{ Row = 18; Nome = "Leonardo" },
    @{ Row = 19; Nome = "Guilherme " },
    @{ Row = 20; Nome = "Wesley" },
    @{ Row = 21; Nome = "Yuri" },
    @{ Row = 22; Nome = "Kira" },
    @{ Row = 23; Nome = "Cleici" },
    @{ Row = 24; Nome = "João Moacir" },
    @{ Row = 25; Nome = "Bruno" },
    @{ Row = 26; Nome = "Elcio" },
    @{ Row = 27; Nome = "Criscia" }
)

foreach ($aluno in $alunos) {
    $r = $aluno.Row

    # Restore the real name (strip "_Modificado")
    $ws.Cells.Item($r, 2).Value = $aluno.Nome

    $faltas = $ws.Cells.Item($r, 3).Value2
    $p1 = $ws.Cells.Item($r, 4).Value2
    $p2 = $ws.Cells.Item($r, 5).Value2
    $p3 = $ws.Cells.Item($r, 6).Value2
    $media = ($p1 + $p2 + $p3) / 3
    $freqPct = ($totalAulas - $faltas) / $totalAulas * 100

    if ($freqPct -lt 75) {
        $situacao = "Reprovado por Falta"
    } elseif ($media -lt 50) {
        $situacao = "Reprovado por Nota"
    } else {
        $situacao = "Aprovado"
    }

    $ws.Cells.Item($r, 7).Value = $situacao
    $ws.Cells.Item($r, 8).Value = 0
}
